# AirHispania Aerodromos.xlsx edit
# "Agregados aerodromos hechos por Luis y actualizada hoja de calculo"
#
# - Marks a batch of aerodromos (rows owned by Luis / Jose) as tested
#   ("si") instead of "no", and highlights them green (matching the
#   existing "tested" rows like A2/A17/...).
# - Row 200 (LEL8) is reassigned from Jose to Tomas.
# - Row 55 (LE58) gets a new note in column I: "En xplane.es".
# - Updates the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AD")

# Green fill used for confirmed/tested rows (matches existing fgColor FF92D050)
$green = 5296274

# Rows whose "Convertido/TESTEADO" status flips from "no" to "si"
# (A column) and get the matching green highlight (A & B columns).
$rows = @(71, 72, 76, 77, 78, 79, 80, 82, 84, 85, 87, 88, 89, 90, 91, 141, 142, 143, 144, 174, 177, 200, 203, 207, 208)

foreach ($r in $rows) {
    $ws.Range("A$r").Value = "si"
    $ws.Range("A$r").Interior.Color = $green
    $ws.Range("B$r").Interior.Color = $green
}

# Row 200 (LEL8) reassigned from Jose to Tomas.
$ws.Range("B200").Value = "Tomas"

# New note for row 55 (LE58).
$ws.Range("I55").Value = "En xplane.es"

# Update the active view/selection to match the edited area.
$ws.Activate()
$ws.Range("I200").Select()
